# Slide 10 ("Исследование задачи"): wording tweak in the "TextBox 4" shape
# (algorithm -> method, Kmeans -> KMeans, and a shorter text box height now
# that the paragraph loses a few characters), plus swapping the vertical
# positions of the elbow-method plot and the k_opt formula picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- "TextBox 4": С помощью алгоритма Локтя строим модель Kmeans, ... ---
$tb = $s.Shapes.Item(3)
$tf = $tb.TextFrame
$tr = $tf.TextRange

# First run: "С помощью алгоритма Локтя строим модель " (40 chars) ->
#            "С помощью метода Локтя строим модель " (37 chars)
$run1 = $tr.Characters(1, 40)
$run1.Text = "С помощью метода Локтя строим модель "

# Second run used to be "Kmeans" right after the first run; re-find it since
# the first run's length just changed (offsets shifted left by 3 chars).
$tr = $tf.TextRange
$run2 = $tr.Find("Kmeans", 0)
$run2.Text = "KMeans"

# The shape uses spAutoFit, and with less text it now sizes to a shorter box.
$tb.Height = 199.52968603937006

# --- "Рисунок 3" (elbow-method inertia plot) moves up ... ---
$elbowPic = $s.Shapes.Item(5)
$elbowPic.Top = 270

# --- ... and "Рисунок 8" (k_opt formula) moves down to take its place ---
$formulaPic = $s.Shapes.Item(6)
$formulaPic.Top = 436.8120582440945
